$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "purpose" column (E2:E25) from "S.GISH" to "fullRNASEQ"
$ws.Range("E2:E25").Value = "fullRNASEQ"

# Reflect the selection change recorded in the saved file (E24:E25)
$ws.Range("E24:E25").Select()
